# Append: 2025-10-30 01:51 JST
# Update the "取得日時" (retrieved timestamp) column for all existing data
# rows to the new run time, and flag row 17's title as featured ("注目").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$oldTimestamp = "2025-10-30 01:21:17"
$newTimestamp = "2025-10-30 01:51:27"

# Data rows run from row 2 through row 18 (row 1 is the header row).
for ($r = 2; $r -le 18; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    if ($cell.Value2 -eq $oldTimestamp) {
        $cell.Value = $newTimestamp
    }
}

# Mark the YouTube music-distribution listing (row 17) as featured.
$titleCell = $ws.Cells.Item(17, 2)
if ($titleCell.Value2 -eq "【急募】YouTubeの音楽配信構築の依頼です") {
    $titleCell.Value = "注目 【急募】YouTubeの音楽配信構築の依頼です"
}
